$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.04741066666666666
$ws.Range("H2").Value = 0.142232
$ws.Range("I2").Value = 0.003188134523263584
$ws.Range("J2").Value = 0.003188134523263585
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.01036033333333333
$ws.Range("N2").Value = 0.031081
$ws.Range("O2").Value = 0.0003369947480386084
$ws.Range("P2").Value = 0.0003369947480386084
$ws.Range("Q2").Value = 0.0004911903102222222
$ws.Range("R2").Value = 0.004420712792
$ws.Range("S2").Value = [double]"1.0743845903804E-06"
$ws.Range("T2").Value = [double]"1.074384590380401E-06"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.04741066666666666
$ws.Range("H3").Value = 0.142232
$ws.Range("I3").Value = 0.003188134523263584
$ws.Range("J3").Value = 0.003188134523263585
$ws.Range("O3").Value = 0.8439700329797517
$ws.Range("P3").Value = 0.8439700329797518
$ws.Range("Q3").Value = 1.230137575527111
$ws.Range("R3").Value = 11.071238179744
$ws.Range("S3").Value = 0.002690689998742652
$ws.Range("T3").Value = 0.002690689998742653
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.04741066666666666
$ws.Range("H4").Value = 0.142232
$ws.Range("I4").Value = 0.003188134523263584
$ws.Range("J4").Value = 0.003188134523263585
$ws.Range("O4").Value = 0.1556929722722096
$ws.Range("P4").Value = 0.1556929722722096
$ws.Range("Q4").Value = 0.2269319619813333
$ws.Range("R4").Value = 2.042387657832
$ws.Range("S4").Value = 0.0004963701399305514
$ws.Range("T4").Value = 0.0004963701399305516
$ws.Range("I5").Value = 0.01595759596384214
$ws.Range("J5").Value = 0.01595759596384214
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.01036033333333333
$ws.Range("N5").Value = 0.031081
$ws.Range("O5").Value = 0.0003369947480386084
$ws.Range("P5").Value = 0.0003369947480386084
$ws.Range("Q5").Value = 0.002458558901666667
$ws.Range("R5").Value = 0.022127030115
$ws.Range("S5").Value = [double]"5.377626031136895E-06"
$ws.Range("T5").Value = [double]"5.377626031136895E-06"
$ws.Range("I6").Value = 0.01595759596384214
$ws.Range("J6").Value = 0.01595759596384214
$ws.Range("O6").Value = 0.8439700329797517
$ws.Range("P6").Value = 0.8439700329797518
$ws.Range("S6").Value = 0.0134677327918814
$ws.Range("T6").Value = 0.0134677327918814
$ws.Range("I7").Value = 0.01595759596384214
$ws.Range("J7").Value = 0.01595759596384214
$ws.Range("O7").Value = 0.1556929722722096
$ws.Range("P7").Value = 0.1556929722722096
$ws.Range("S7").Value = 0.002484485545929597
$ws.Range("T7").Value = 0.002484485545929598
$ws.Range("I8").Value = 0.9808542695128942
$ws.Range("J8").Value = 0.9808542695128943
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.01036033333333333
$ws.Range("N8").Value = 0.031081
$ws.Range("O8").Value = 0.0003369947480386084
$ws.Range("P8").Value = 0.0003369947480386084
$ws.Range("Q8").Value = 0.1511185018728889
$ws.Range("R8").Value = 1.360066516856
$ws.Range("S8").Value = 0.0003305427374170911
$ws.Range("T8").Value = 0.0003305427374170911
$ws.Range("I9").Value = 0.9808542695128942
$ws.Range("J9").Value = 0.9808542695128943
$ws.Range("O9").Value = 0.8439700329797517
$ws.Range("P9").Value = 0.8439700329797518
$ws.Range("S9").Value = 0.8278116101891276
$ws.Range("T9").Value = 0.8278116101891279
$ws.Range("I10").Value = 0.9808542695128942
$ws.Range("J10").Value = 0.9808542695128943
$ws.Range("O10").Value = 0.1556929722722096
$ws.Range("P10").Value = 0.1556929722722096
$ws.Range("S10").Value = 0.1527121165863495
$ws.Range("T10").Value = 0.1527121165863495
